$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 from 3 to 5
$ws.Range("B2").Value = 5

# Update A3 from 2 to 1
$ws.Range("A3").Value = 1

# Remove row 4 entirely (A4=1, B4=2), shrinking the used range to A1:B3
$ws.Rows("4:4").Delete()
